# Update the "processed_classification_result" sheet:
#   - rename the "humoment" feature group to "shape" in the header row
#   - replace the stale per-classifier metric values (rows 4-7) with the
#     numbers from the latest classification run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the `humoment` feature-group header cells to `shape`
# (base header + the three merged combination headers that include it).
$ws.Range("J1").Value = "shape"
$ws.Range("R1").Value = "texture-shape"
$ws.Range("V1").Value = "color-shape"
$ws.Range("Z1").Value = "texture-color-shape"

# Refresh the classification metrics (f1 / recall / precision / accuracy for
# each of the 7 feature-group columns) for every classifier row.
# Row 4
$ws.Cells.Item(4, 2).Value = 0.5932486733548871
$ws.Cells.Item(4, 3).Value = 0.622
$ws.Cells.Item(4, 4).Value = 0.580538685613357
$ws.Cells.Item(4, 5).Value = 0.604
$ws.Cells.Item(4, 6).Value = 0.6084982694032386
$ws.Cells.Item(4, 7).Value = 0.6499999999999999
$ws.Cells.Item(4, 8).Value = 0.6151592825421899
$ws.Cells.Item(4, 9).Value = 0.6165
$ws.Cells.Item(4, 10).Value = 0.6698229371633928
$ws.Cells.Item(4, 11).Value = 0.9799999999999999
$ws.Cells.Item(4, 12).Value = 0.5089641738515409
$ws.Cells.Item(4, 13).Value = 0.517
$ws.Cells.Item(4, 14).Value = 0.6436996983212654
$ws.Cells.Item(4, 15).Value = 0.6599999999999999
$ws.Cells.Item(4, 16).Value = 0.6641647457149111
$ws.Cells.Item(4, 17).Value = 0.6570000000000001
$ws.Cells.Item(4, 18).Value = 0.5900637650000664
$ws.Cells.Item(4, 19).Value = 0.6199999999999999
$ws.Cells.Item(4, 20).Value = 0.57565129568173
$ws.Cells.Item(4, 21).Value = 0.6004999999999999
$ws.Cells.Item(4, 22).Value = 0.601833407479299
$ws.Cells.Item(4, 23).Value = 0.6449999999999999
$ws.Cells.Item(4, 24).Value = 0.6072180444423341
$ws.Cells.Item(4, 25).Value = 0.61
$ws.Cells.Item(4, 26).Value = 0.6420919044424304
$ws.Cells.Item(4, 27).Value = 0.659
$ws.Cells.Item(4, 28).Value = 0.6623320251403092
$ws.Cells.Item(4, 29).Value = 0.655
# Row 5
$ws.Cells.Item(5, 2).Value = 0.5963164773123009
$ws.Cells.Item(5, 3).Value = 0.624
$ws.Cells.Item(5, 4).Value = 0.5825920123720698
$ws.Cells.Item(5, 5).Value = 0.6055
$ws.Cells.Item(5, 6).Value = 0.7035724931051273
$ws.Cells.Item(5, 7).Value = 0.765
$ws.Cells.Item(5, 8).Value = 0.6674766636113154
$ws.Cells.Item(5, 9).Value = 0.675
$ws.Cells.Item(5, 10).Value = 0.6525963824887014
$ws.Cells.Item(5, 11).Value = 0.857
$ws.Cells.Item(5, 12).Value = 0.5320614704031437
$ws.Cells.Item(5, 13).Value = 0.5485
$ws.Cells.Item(5, 14).Value = 0.6389404811190171
$ws.Cells.Item(5, 15).Value = 0.652
$ws.Cells.Item(5, 16).Value = 0.6481075370780268
$ws.Cells.Item(5, 17).Value = 0.6519999999999999
$ws.Cells.Item(5, 18).Value = 0.5964912338594319
$ws.Cells.Item(5, 19).Value = 0.6239999999999999
$ws.Cells.Item(5, 20).Value = 0.5829737464980567
$ws.Cells.Item(5, 21).Value = 0.6054999999999999
$ws.Cells.Item(5, 22).Value = 0.6832539814881928
$ws.Cells.Item(5, 23).Value = 0.7470000000000001
$ws.Cells.Item(5, 24).Value = 0.6493748164999953
$ws.Cells.Item(5, 25).Value = 0.6525000000000001
$ws.Cells.Item(5, 26).Value = 0.6391019338593521
$ws.Cells.Item(5, 27).Value = 0.6540000000000001
$ws.Cells.Item(5, 28).Value = 0.6471785098518685
$ws.Cells.Item(5, 29).Value = 0.652
# Row 6
$ws.Cells.Item(6, 2).Value = 0.6036519292314397
$ws.Cells.Item(6, 3).Value = 0.635
$ws.Cells.Item(6, 4).Value = 0.5912403920193882
$ws.Cells.Item(6, 5).Value = 0.6134999999999999
$ws.Cells.Item(6, 6).Value = 0.7065737934675792
$ws.Cells.Item(6, 7).Value = 0.74
$ws.Cells.Item(6, 8).Value = 0.6930043549100147
$ws.Cells.Item(6, 9).Value = 0.6944999999999999
$ws.Cells.Item(6, 10).Value = 0.6580562848705325
$ws.Cells.Item(6, 11).Value = 0.885
$ws.Cells.Item(6, 12).Value = 0.5272875985295125
$ws.Cells.Item(6, 13).Value = 0.543
$ws.Cells.Item(6, 14).Value = 0.6648943608648168
$ws.Cells.Item(6, 15).Value = 0.6699999999999999
$ws.Cells.Item(6, 16).Value = 0.6811693287624798
$ws.Cells.Item(6, 17).Value = 0.68
$ws.Cells.Item(6, 18).Value = 0.6056245702751673
$ws.Cells.Item(6, 19).Value = 0.6380000000000001
$ws.Cells.Item(6, 20).Value = 0.5927528505075299
$ws.Cells.Item(6, 21).Value = 0.6165
$ws.Cells.Item(6, 22).Value = 0.6635843691012439
$ws.Cells.Item(6, 23).Value = 0.706
$ws.Cells.Item(6, 24).Value = 0.6646969646848756
$ws.Cells.Item(6, 25).Value = 0.653
$ws.Cells.Item(6, 26).Value = 0.655087619539876
$ws.Cells.Item(6, 27).Value = 0.657
$ws.Cells.Item(6, 28).Value = 0.6761167703753674
$ws.Cells.Item(6, 29).Value = 0.673
# Row 7
$ws.Cells.Item(7, 2).Value = 0.4765152655222479
$ws.Cells.Item(7, 3).Value = 0.488
$ws.Cells.Item(7, 4).Value = 0.473380765413255
$ws.Cells.Item(7, 5).Value = 0.4784999999999999
$ws.Cells.Item(7, 6).Value = 0.4839927403138217
$ws.Cells.Item(7, 7).Value = 0.499
$ws.Cells.Item(7, 8).Value = 0.4805376887905241
$ws.Cells.Item(7, 9).Value = 0.4834999999999999
$ws.Cells.Item(7, 10).Value = 0.6183677339517268
$ws.Cells.Item(7, 11).Value = 0.8370000000000001
$ws.Cells.Item(7, 12).Value = 0.4996152807756056
$ws.Cells.Item(7, 13).Value = 0.505
$ws.Cells.Item(7, 14).Value = 0.5178702970339001
$ws.Cells.Item(7, 15).Value = 0.541
$ws.Cells.Item(7, 16).Value = 0.511124862071197
$ws.Cells.Item(7, 17).Value = 0.5170000000000001
$ws.Cells.Item(7, 18).Value = 0.4907741497664189
$ws.Cells.Item(7, 19).Value = 0.506
$ws.Cells.Item(7, 20).Value = 0.484064698609001
$ws.Cells.Item(7, 21).Value = 0.487
$ws.Cells.Item(7, 22).Value = 0.5347117507663307
$ws.Cells.Item(7, 23).Value = 0.5559999999999999
$ws.Cells.Item(7, 24).Value = 0.5410911672502035
$ws.Cells.Item(7, 25).Value = 0.5475000000000001
$ws.Cells.Item(7, 26).Value = 0.5260243126713876
$ws.Cells.Item(7, 27).Value = 0.549
$ws.Cells.Item(7, 28).Value = 0.5198667432613931
$ws.Cells.Item(7, 29).Value = 0.526
